$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with new TPM-derived values
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna2"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.015377
$ws.Range("H2").Value = 6.046131
$ws.Range("I2").Value = 0.7554960962715589
$ws.Range("J2").Value = 0.7554960962715588
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.081040666666667
$ws.Range("N2").Value = 24.243122
$ws.Range("O2").Value = 0.4661250698616886
$ws.Range("P2").Value = 0.4661250698616886
$ws.Range("Q2").Value = 16.28634349566467
$ws.Range("R2").Value = 146.577091460982
$ws.Range("S2").Value = 0.3521556706548134
$ws.Range("T2").Value = 0.3521556706548133

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna2"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.015377
$ws.Range("H3").Value = 6.046131
$ws.Range("I3").Value = 0.7554960962715589
$ws.Range("J3").Value = 0.7554960962715588
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.000300666666668
$ws.Range("N3").Value = 21.000902
$ws.Range("O3").Value = 0.4037865631294714
$ws.Range("P3").Value = 0.4037865631294715
$ws.Range("Q3").Value = 14.10824495668467
$ws.Range("R3").Value = 126.974204610162
$ws.Range("S3").Value = 0.305059172171225
$ws.Range("T3").Value = 0.305059172171225

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna2"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.015377
$ws.Range("H4").Value = 6.046131
$ws.Range("I4").Value = 0.7554960962715589
$ws.Range("J4").Value = 0.7554960962715588
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.255294666666666
$ws.Range("N4").Value = 6.765884
$ws.Range("O4").Value = 0.1300883670088399
$ws.Range("P4").Value = 0.1300883670088399
$ws.Range("Q4").Value = 4.545268999422666
$ws.Range("R4").Value = 40.907420994804
$ws.Range("S4").Value = 0.0982812534455204
$ws.Range("T4").Value = 0.09828125344552038

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna2"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6522436666666667
$ws.Range("H5").Value = 1.956731
$ws.Range("I5").Value = 0.2445039037284412
$ws.Range("J5").Value = 0.2445039037284411
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.081040666666667
$ws.Range("N5").Value = 24.243122
$ws.Range("O5").Value = 0.4661250698616886
$ws.Range("P5").Value = 0.4661250698616886
$ws.Range("Q5").Value = 5.270807594909111
$ws.Range("R5").Value = 47.437268354182
$ws.Range("S5").Value = 0.1139693992068752
$ws.Range("T5").Value = 0.1139693992068752

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna2"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6522436666666667
$ws.Range("H6").Value = 1.956731
$ws.Range("I6").Value = 0.2445039037284412
$ws.Range("J6").Value = 0.2445039037284411
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.000300666666668
$ws.Range("N6").Value = 21.000902
$ws.Range("O6").Value = 0.4037865631294714
$ws.Range("P6").Value = 0.4037865631294715
$ws.Range("Q6").Value = 4.565901774595779
$ws.Range("R6").Value = 41.093115971362
$ws.Range("S6").Value = 0.09872739095824641
$ws.Range("T6").Value = 0.09872739095824641

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna2"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6522436666666667
$ws.Range("H7").Value = 1.956731
$ws.Range("I7").Value = 0.2445039037284412
$ws.Range("J7").Value = 0.2445039037284411
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.255294666666666
$ws.Range("N7").Value = 6.765884
$ws.Range("O7").Value = 0.1300883670088399
$ws.Range("P7").Value = 0.1300883670088399
$ws.Range("Q7").Value = 1.471001662800444
$ws.Range("R7").Value = 13.239014965204
$ws.Range("S7").Value = 0.03180711356331951
$ws.Range("T7").Value = 0.03180711356331951

# Remove now-obsolete rows 8-10 (data now fits in rows 2-7)
$ws.Rows("8:10").Delete()
